$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 241: Jin Bu-Yeon quote about the ice stone
$ws.Range("A241").Value = "No one can truly own the ice stone. You only use its power. Pull and get pushed. Earn one, and lose the other. That is the burden you will have to carry."
$ws.Range("B241").Value = "Jin Bu-Yeon"
$ws.Range("C241").Value = 1
$ws.Range("D241").Value = 18

# Row 242: Jang Uk quote about embarrassment
$ws.Range("A242").Value = "I can always say such things﻿ if I endure the embarrassment. But I can never turn back time even if I regret something. Master, endure the embarrassment and be honest with yourself. That way, you will not regret anything."
$ws.Range("B242").Value = "Jang Uk"
$ws.Range("C242").Value = 1
$ws.Range("D242").Value = 18

# Row 243: Jang Uk & Naksu exchange
$ws.Range("B243").Value = "Jang Uk & Naksu"
$ws.Range("A243").Value = "(Jang Uk) Why did you hide it? (Naksu) Because…I wanted you to die. (Jang Uk) Then…why did you come for me? (Naksu) Because…I thought you would die."
$ws.Range("C243").Value = 1
$ws.Range("D243").Value = 18

$ws.Range("A241").Select()
